# The source workbook is a VEDA "~TFM_INS" table. Its TimeSlice column
# (column B) contained the value "DAYNITE" in every data row; the edit
# changes that value to "ANNUAL" (a shared-string in-place update,
# equivalent to a find/replace over the used range), and leaves the
# worksheet selection on B10 (the last populated TimeSlice cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        if ($cell.Value() -eq "DAYNITE") {
            $cell.Value = "ANNUAL"
        }
    }
}

# Move / record the active selection, matching the saved view state.
$ws.Range("B10").Select()
